# Update "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps on the zh-cn and de-de report sheets, as part of
# regenerating the handback status report.

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")

# Row 2
$wsZh.Range("E2").Value = "2016-03-18 03:27:07"
$wsZh.Range("H2").Value = "2016-03-18 03:27:50"

# Row 3
$wsZh.Range("E3").Value = "2016-03-18 03:27:07"
$wsZh.Range("H3").Value = "2016-03-18 03:27:50"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")

# Row 2
$wsDe.Range("E2").Value = "2016-03-18 03:27:15"
$wsDe.Range("H2").Value = "2016-03-18 03:28:03"

# Row 3
$wsDe.Range("E3").Value = "2016-03-18 03:27:15"
$wsDe.Range("H3").Value = "2016-03-18 03:28:03"
